$d = $word.ActiveDocument
$failCount = 0

if (-not $d.Content.Find.Execute("55+7=", $true, $true, $false, $false, $false, $true, 1, $false, "74+7=", 1)) { $failCount++ }
if (-not $d.Content.Find.Execute("99-91=", $true, $true, $false, $false, $false, $true, 1, $false, "81-78=", 1)) { $failCount++ }
if (-not $d.Content.Find.Execute("81-28=", $true, $true, $false, $false, $false, $true, 1, $false, "67-26=", 1)) { $failCount++ }
if (-not $d.Content.Find.Execute("43+15=", $true, $true, $false, $false, $false, $true, 1, $false, "84-17=", 1)) { $failCount++ }
if (-not $d.Content.Find.Execute("31+52=", $true, $true, $false, $false, $false, $true, 1, $false, "99-16=", 1)) { $failCount++ }
if (-not $d.Content.Find.Execute("40+39=", $true, $true, $false, $false, $false, $true, 1, $false, "9+13=", 1)) { $failCount++ }
if (-not $d.Content.Find.Execute("16+57=", $true, $true, $false, $false, $false, $true, 1, $false, "51-32=", 1)) { $failCount++ }
if (-not $d.Content.Find.Execute("85-59=", $true, $true, $false, $false, $false, $true, 1, $false, "36-35=", 1)) { $failCount++ }
if (-not $d.Content.Find.Execute("90-76=", $true, $true, $false, $false, $false, $true, 1, $false, "80-34=", 1)) { $failCount++ }
if (-not $d.Content.Find.Execute("96-37=", $true, $true, $false, $false, $false, $true, 1, $false, "90-44=", 1)) { $failCount++ }
if (-not $d.Content.Find.Execute("77+6=", $true, $true, $false, $false, $false, $true, 1, $false, "28+45=", 1)) { $failCount++ }
if (-not $d.Content.Find.Execute("4+63=", $true, $true, $false, $false, $false, $true, 1, $false, "43+35=", 1)) { $failCount++ }
if (-not $d.Content.Find.Execute("81-5=", $true, $true, $false, $false, $false, $true, 1, $false, "3+63=", 1)) { $failCount++ }
if (-not $d.Content.Find.Execute("0+52=", $true, $true, $false, $false, $false, $true, 1, $false, "54+23=", 1)) { $failCount++ }
if (-not $d.Content.Find.Execute("44-3=", $true, $true, $false, $false, $false, $true, 1, $false, "96-6=", 1)) { $failCount++ }
if (-not $d.Content.Find.Execute("0+83=", $true, $true, $false, $false, $false, $true, 1, $false, "17-4=", 1)) { $failCount++ }
if (-not $d.Content.Find.Execute("14+45=", $true, $true, $false, $false, $false, $true, 1, $false, "15+4=", 1)) { $failCount++ }
if (-not $d.Content.Find.Execute("31+48=", $true, $true, $false, $false, $false, $true, 1, $false, "53-9=", 1)) { $failCount++ }
if (-not $d.Content.Find.Execute("28+21=", $true, $true, $false, $false, $false, $true, 1, $false, "40-3=", 1)) { $failCount++ }
if (-not $d.Content.Find.Execute("83-78=", $true, $true, $false, $false, $false, $true, 1, $false, "63+18=", 1)) { $failCount++ }
if (-not $d.Content.Find.Execute("11+63=", $true, $true, $false, $false, $false, $true, 1, $false, "32-2=", 1)) { $failCount++ }
if (-not $d.Content.Find.Execute("10-1=", $true, $true, $false, $false, $false, $true, 1, $false, "84-2=", 1)) { $failCount++ }
if (-not $d.Content.Find.Execute("83+3=", $true, $true, $false, $false, $false, $true, 1, $false, "22+22=", 1)) { $failCount++ }
if (-not $d.Content.Find.Execute("56+23=", $true, $true, $false, $false, $false, $true, 1, $false, "1+49=", 1)) { $failCount++ }
if (-not $d.Content.Find.Execute("80-31=", $true, $true, $false, $false, $false, $true, 1, $false, "29+39=", 1)) { $failCount++ }
if (-not $d.Content.Find.Execute("95-57=", $true, $true, $false, $false, $false, $true, 1, $false, "96-77=", 1)) { $failCount++ }
if (-not $d.Content.Find.Execute("32+45=", $true, $true, $false, $false, $false, $true, 1, $false, "45-7=", 1)) { $failCount++ }
if (-not $d.Content.Find.Execute("90-40=", $true, $true, $false, $false, $false, $true, 1, $false, "21+63=", 1)) { $failCount++ }
if (-not $d.Content.Find.Execute("70-66=", $true, $true, $false, $false, $false, $true, 1, $false, "13+26=", 1)) { $failCount++ }
if (-not $d.Content.Find.Execute("25+18=", $true, $true, $false, $false, $false, $true, 1, $false, "32-27=", 1)) { $failCount++ }
if (-not $d.Content.Find.Execute("29+49=", $true, $true, $false, $false, $false, $true, 1, $false, "42-36=", 1)) { $failCount++ }
if (-not $d.Content.Find.Execute("2+14=", $true, $true, $false, $false, $false, $true, 1, $false, "33+38=", 1)) { $failCount++ }
if (-not $d.Content.Find.Execute("52+12=", $true, $true, $false, $false, $false, $true, 1, $false, "5+29=", 1)) { $failCount++ }
if (-not $d.Content.Find.Execute("83-39=", $true, $true, $false, $false, $false, $true, 1, $false, "25-8=", 1)) { $failCount++ }
if (-not $d.Content.Find.Execute("27+44=", $true, $true, $false, $false, $false, $true, 1, $false, "8+42=", 1)) { $failCount++ }
if (-not $d.Content.Find.Execute("59-5=", $true, $true, $false, $false, $false, $true, 1, $false, "78-70=", 1)) { $failCount++ }
if (-not $d.Content.Find.Execute("87-12=", $true, $true, $false, $false, $false, $true, 1, $false, "69+12=", 1)) { $failCount++ }
if (-not $d.Content.Find.Execute("98-40=", $true, $true, $false, $false, $false, $true, 1, $false, "26+9=", 1)) { $failCount++ }
if (-not $d.Content.Find.Execute("63-0=", $true, $true, $false, $false, $false, $true, 1, $false, "79-48=", 1)) { $failCount++ }
if (-not $d.Content.Find.Execute("12+44=", $true, $true, $false, $false, $false, $true, 1, $false, "72-2=", 1)) { $failCount++ }
if (-not $d.Content.Find.Execute("99-18=", $true, $true, $false, $false, $false, $true, 1, $false, "27+71=", 1)) { $failCount++ }
if (-not $d.Content.Find.Execute("0+75=", $true, $true, $false, $false, $false, $true, 1, $false, "28-6=", 1)) { $failCount++ }
if (-not $d.Content.Find.Execute("53+15=", $true, $true, $false, $false, $false, $true, 1, $false, "2+57=", 1)) { $failCount++ }
if (-not $d.Content.Find.Execute("70-64=", $true, $true, $false, $false, $false, $true, 1, $false, "55+16=", 1)) { $failCount++ }
if (-not $d.Content.Find.Execute("56-8=", $true, $true, $false, $false, $false, $true, 1, $false, "76-20=", 1)) { $failCount++ }
if (-not $d.Content.Find.Execute("75-25=", $true, $true, $false, $false, $false, $true, 1, $false, "33+0=", 1)) { $failCount++ }
if (-not $d.Content.Find.Execute("74-39=", $true, $true, $false, $false, $false, $true, 1, $false, "61-13=", 1)) { $failCount++ }
if (-not $d.Content.Find.Execute("3+59=", $true, $true, $false, $false, $false, $true, 1, $false, "21+51=", 1)) { $failCount++ }
if (-not $d.Content.Find.Execute("98-77=", $true, $true, $false, $false, $false, $true, 1, $false, "54-17=", 1)) { $failCount++ }
if (-not $d.Content.Find.Execute("81-56=", $true, $true, $false, $false, $false, $true, 1, $false, "40+49=", 1)) { $failCount++ }
if (-not $d.Content.Find.Execute("32-12=", $true, $true, $false, $false, $false, $true, 1, $false, "34+6=", 1)) { $failCount++ }
if (-not $d.Content.Find.Execute("33-21=", $true, $true, $false, $false, $false, $true, 1, $false, "45+32=", 1)) { $failCount++ }
if (-not $d.Content.Find.Execute("24+57=", $true, $true, $false, $false, $false, $true, 1, $false, "51-26=", 1)) { $failCount++ }
if (-not $d.Content.Find.Execute("72-66=", $true, $true, $false, $false, $false, $true, 1, $false, "14+25=", 1)) { $failCount++ }
if (-not $d.Content.Find.Execute("93-73=", $true, $true, $false, $false, $false, $true, 1, $false, "2+48=", 1)) { $failCount++ }
if (-not $d.Content.Find.Execute("57+3=", $true, $true, $false, $false, $false, $true, 1, $false, "72-66=", 1)) { $failCount++ }
if (-not $d.Content.Find.Execute("54+15=", $true, $true, $false, $false, $false, $true, 1, $false, "32-7=", 1)) { $failCount++ }
if (-not $d.Content.Find.Execute("67-47=", $true, $true, $false, $false, $false, $true, 1, $false, "1+46=", 1)) { $failCount++ }
if (-not $d.Content.Find.Execute("80-69=", $true, $true, $false, $false, $false, $true, 1, $false, "90-68=", 1)) { $failCount++ }
if (-not $d.Content.Find.Execute("20-18=", $true, $true, $false, $false, $false, $true, 1, $false, "33-17=", 1)) { $failCount++ }
if (-not $d.Content.Find.Execute("81+12=", $true, $true, $false, $false, $false, $true, 1, $false, "54+7=", 1)) { $failCount++ }
if (-not $d.Content.Find.Execute("64+4=", $true, $true, $false, $false, $false, $true, 1, $false, "73+24=", 1)) { $failCount++ }
if (-not $d.Content.Find.Execute("60+7=", $true, $true, $false, $false, $false, $true, 1, $false, "76+5=", 1)) { $failCount++ }
if (-not $d.Content.Find.Execute("42-20=", $true, $true, $false, $false, $false, $true, 1, $false, "39+36=", 1)) { $failCount++ }
if (-not $d.Content.Find.Execute("18+38=", $true, $true, $false, $false, $false, $true, 1, $false, "50+3=", 1)) { $failCount++ }
if (-not $d.Content.Find.Execute("8-6=", $true, $true, $false, $false, $false, $true, 1, $false, "84-42=", 1)) { $failCount++ }
if (-not $d.Content.Find.Execute("79+8=", $true, $true, $false, $false, $false, $true, 1, $false, "16+28=", 1)) { $failCount++ }
if (-not $d.Content.Find.Execute("4+39=", $true, $true, $false, $false, $false, $true, 1, $false, "26-23=", 1)) { $failCount++ }
if (-not $d.Content.Find.Execute("29+48=", $true, $true, $false, $false, $false, $true, 1, $false, "86-7=", 1)) { $failCount++ }
if (-not $d.Content.Find.Execute("99-71=", $true, $true, $false, $false, $false, $true, 1, $false, "16-15=", 1)) { $failCount++ }
if (-not $d.Content.Find.Execute("97-84=", $true, $true, $false, $false, $false, $true, 1, $false, "2+67=", 1)) { $failCount++ }
if (-not $d.Content.Find.Execute("56+29=", $true, $true, $false, $false, $false, $true, 1, $false, "61-26=", 1)) { $failCount++ }
if (-not $d.Content.Find.Execute("25+31=", $true, $true, $false, $false, $false, $true, 1, $false, "29+32=", 1)) { $failCount++ }
if (-not $d.Content.Find.Execute("88+8=", $true, $true, $false, $false, $false, $true, 1, $false, "63-13=", 1)) { $failCount++ }
if (-not $d.Content.Find.Execute("52+12=", $true, $true, $false, $false, $false, $true, 1, $false, "82-39=", 1)) { $failCount++ }
if (-not $d.Content.Find.Execute("11+14=", $true, $true, $false, $false, $false, $true, 1, $false, "36-5=", 1)) { $failCount++ }
if (-not $d.Content.Find.Execute("62-24=", $true, $true, $false, $false, $false, $true, 1, $false, "50+36=", 1)) { $failCount++ }
if (-not $d.Content.Find.Execute("5-1=", $true, $true, $false, $false, $false, $true, 1, $false, "3-0=", 1)) { $failCount++ }
if (-not $d.Content.Find.Execute("71-61=", $true, $true, $false, $false, $false, $true, 1, $false, "58+10=", 1)) { $failCount++ }
if (-not $d.Content.Find.Execute("18+58=", $true, $true, $false, $false, $false, $true, 1, $false, "93-4=", 1)) { $failCount++ }
if (-not $d.Content.Find.Execute("93-52=", $true, $true, $false, $false, $false, $true, 1, $false, "8+81=", 1)) { $failCount++ }
if (-not $d.Content.Find.Execute("45+40=", $true, $true, $false, $false, $false, $true, 1, $false, "0+32=", 1)) { $failCount++ }
if (-not $d.Content.Find.Execute("27+66=", $true, $true, $false, $false, $false, $true, 1, $false, "28+42=", 1)) { $failCount++ }
if (-not $d.Content.Find.Execute("98-74=", $true, $true, $false, $false, $false, $true, 1, $false, "12+15=", 1)) { $failCount++ }
if (-not $d.Content.Find.Execute("70-63=", $true, $true, $false, $false, $false, $true, 1, $false, "77+22=", 1)) { $failCount++ }
if (-not $d.Content.Find.Execute("32+9=", $true, $true, $false, $false, $false, $true, 1, $false, "32+29=", 1)) { $failCount++ }
if (-not $d.Content.Find.Execute("50-13=", $true, $true, $false, $false, $false, $true, 1, $false, "65-57=", 1)) { $failCount++ }
if (-not $d.Content.Find.Execute("23-10=", $true, $true, $false, $false, $false, $true, 1, $false, "1+19=", 1)) { $failCount++ }
if (-not $d.Content.Find.Execute("0+16=", $true, $true, $false, $false, $false, $true, 1, $false, "26+24=", 1)) { $failCount++ }
if (-not $d.Content.Find.Execute("89+2=", $true, $true, $false, $false, $false, $true, 1, $false, "11+31=", 1)) { $failCount++ }
if (-not $d.Content.Find.Execute("12+21=", $true, $true, $false, $false, $false, $true, 1, $false, "95-32=", 1)) { $failCount++ }
if (-not $d.Content.Find.Execute("59-15=", $true, $true, $false, $false, $false, $true, 1, $false, "66-24=", 1)) { $failCount++ }
if (-not $d.Content.Find.Execute("32-3=", $true, $true, $false, $false, $false, $true, 1, $false, "61-47=", 1)) { $failCount++ }
if (-not $d.Content.Find.Execute("21-17=", $true, $true, $false, $false, $false, $true, 1, $false, "39+45=", 1)) { $failCount++ }
if (-not $d.Content.Find.Execute("41+31=", $true, $true, $false, $false, $false, $true, 1, $false, "62-54=", 1)) { $failCount++ }
if (-not $d.Content.Find.Execute("95-58=", $true, $true, $false, $false, $false, $true, 1, $false, "52-42=", 1)) { $failCount++ }
if (-not $d.Content.Find.Execute("63-18=", $true, $true, $false, $false, $false, $true, 1, $false, "30+28=", 1)) { $failCount++ }
if (-not $d.Content.Find.Execute("12+80=", $true, $true, $false, $false, $false, $true, 1, $false, "13+25=", 1)) { $failCount++ }
if (-not $d.Content.Find.Execute("13+39=", $true, $true, $false, $false, $false, $true, 1, $false, "97-49=", 1)) { $failCount++ }
if (-not $d.Content.Find.Execute("52-32=", $true, $true, $false, $false, $false, $true, 1, $false, "38+3=", 1)) { $failCount++ }

if ($failCount -gt 0) { Write-Output "Warning: $failCount replacement(s) failed" } else { Write-Output "All replacements applied successfully" }
